# Add "2022-Q1" worksheet (positioned right before the "总计" summary sheet)
# and populate it with the new quarter's fund-holdings data, then prepend
# a matching summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$totalWsBefore = $wb.Worksheets.Item("总计")
$sourceWs = $wb.Worksheets.Item("2021-Q3")

# --- 1. Create the new "2022-Q1" sheet just before "总计" -------------------
$newWs = $wb.Worksheets.Add($totalWsBefore)
$newWs.Name = "2022-Q1"

# NOTE: passing a worksheet object into Add() rebinds that variable to the
# freshly created sheet in this engine, so re-fetch a *fresh* reference to
# the "总计" sheet by name now that it has been pushed one slot further on.
$totalWs = $wb.Worksheets.Item("总计")

# Copy layout/formatting (header row + 3 data rows) from an existing
# quarterly sheet so fonts/borders/styles match the rest of the workbook.
$sourceWs.Range("A1:H4").Copy($newWs.Range("A1"))

# --- header row ---
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Columns B, D, E, F, G hold text-like values (fund codes / percentages
# formatted as text) in the source data, so force Text format before
# assigning to avoid Excel auto-converting them to numbers (which would
# strip leading zeros / alter precision).
$newWs.Range("B2:B4").NumberFormat = "@"
$newWs.Range("D2:G4").NumberFormat = "@"

# --- row 2 : 070019 嘉实价值优势混合 ---
$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "070019"
$newWs.Range("C2").Value = "嘉实价值优势混合"
$newWs.Range("D2").Value = "27.51"
$newWs.Range("E2").Value = "93.14"
$newWs.Range("F2").Value = "5.13"
$newWs.Range("G2").Value = "1.4113"
$newWs.Range("H2").Value = 9

# --- row 3 : 001075 宝盈转型动力灵活配置混合 ---
$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "001075"
$newWs.Range("C3").Value = "宝盈转型动力灵活配置混合"
$newWs.Range("D3").Value = "5.13"
$newWs.Range("E3").Value = "86.64"
$newWs.Range("F3").Value = "4.48"
$newWs.Range("G3").Value = "0.2298"
$newWs.Range("H3").Value = 1

# --- row 4 : 010676 光大保德信新机遇混合 ---
$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "010676"
$newWs.Range("C4").Value = "光大保德信新机遇混合"
$newWs.Range("D4").Value = "4.20"
$newWs.Range("E4").Value = "85.57"
$newWs.Range("F4").Value = "5.47"
$newWs.Range("G4").Value = "0.2297"
$newWs.Range("H4").Value = 3

# --- 2. Insert a new top data row in "总计" for 2022-Q1 --------------------
$totalWs.Rows(2).Insert()

# The freshly inserted row inherits the (bold/centered) formatting of the
# header row above it. Re-apply the correct plain-data-row formatting by
# copying it from the row directly below, which already carries it
# (index column keeps its centered style, the rest stay unstyled).
$totalWs.Range("A3:D3").Copy($totalWs.Range("A2"))

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 3
$totalWs.Range("D2").Value = 1.87

# Renumber the index column for the rows that got shifted down.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
